$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "cod_api" column header in G2, matching style of the other header cells
$ws.Range("G2").Value = "cod_api"
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)  # xlPasteFormats

# Fill G3:G8 with 0, matching style of the neighbouring numeric columns (e.g. E column)
$ws.Range("G3:G8").Value = 0
$ws.Range("E3:E8").Copy()
$ws.Range("G3:G8").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the selection, as Excel does after editing the new column
$ws.Range("G3").Select()
